$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$schemaJson = @'
{
    "$schema": "http://json-schema.org/draft-07/schema#",
    "description": "MOSIP Sample identity",
    "additionalProperties": false,
    "title": "MOSIP identity",
    "type": "object",
    "definitions": {
        "simpleType": {
            "uniqueItems": true,
            "additionalItems": false,
            "type": "array",
            "items": {
                "additionalProperties": false,
                "type": "object",
                "required": [
                    "language",
                    "value"
                ],
                "properties": {
                    "language": {
                        "type": "string"
                    },
                    "value": {
                        "type": "string"
                    }
                }
            }
        },
        "documentType": {
            "additionalProperties": false,
            "type": "object",
            "properties": {
                "format": {
                    "type": "string"
                },
                "type": {
                    "type": "string"
                },
                "value": {
                    "type": "string"
                },
                "refNumber": {
                    "type": [
                        "string",
                        "null"
                    ]
                }
            }
        },
        "biometricsType": {
            "additionalProperties": false,
            "type": "object",
            "properties": {
                "format": {
                    "type": "string"
                },
                "version": {
                    "type": "number",
                    "minimum": 0
                },
                "value": {
                    "type": "string"
                }
            }
        }
    },
    "properties": {
        "identity": {
            "additionalProperties": false,
            "type": "object",
            "required": [
                "IDSchemaVersion",
                "firstName",
                "lastName",
                "dateOfBirth",
                "gender",
                "Region",
                "Departement",
                "ChefLieu",
                "SousPrefectures",
                "countryOfCitizenship",
                "individualBiometrics",
                "residenceStatus"
            ],
            "properties": {
                "proofOfCnamEnrollment": {
                    "bioAttributes": [],
                    "fieldCategory": "pvt",
                    "format": "none",
                    "fieldType": "default",
                    "$ref": "#/definitions/documentType"
                },
                "gender": {
                    "bioAttributes": [],
                    "fieldCategory": "pvt",
                    "format": "",
                    "fieldType": "default",
                    "$ref": "#/definitions/simpleType"
                },
                "Region": {
                    "bioAttributes": [],
                    "fieldCategory": "pvt",
                    "format": "none",
                    "fieldType": "default",
                    "$ref": "#/definitions/simpleType"
                },
                "Departement": {
                    "bioAttributes": [],
                    "fieldCategory": "pvt",
                    "format": "none",
                    "fieldType": "default",
                    "$ref": "#/definitions/simpleType"
                },
                "ChefLieu": {
                    "bioAttributes": [],
                    "fieldCategory": "pvt",
                    "format": "none",
                    "fieldType": "default",
                    "$ref": "#/definitions/simpleType"
                },
                "SousPrefectures": {
                    "bioAttributes": [],
                    "fieldCategory": "pvt",
                    "format": "none",
                    "fieldType": "default",
                    "$ref": "#/definitions/simpleType"
                },
                "Commune": {
                    "bioAttributes": [],
                    "fieldCategory": "pvt",
                    "format": "none",
                    "fieldType": "default",
                    "$ref": "#/definitions/simpleType"
                },
                "countryOfCitizenship": {
                    "bioAttributes": [],
                    "fieldCategory": "pvt",
                    "format": "none",
                    "fieldType": "default",
                    "$ref": "#/definitions/simpleType"
                },
                "residenceStatus": {
                    "bioAttributes": [],
                    "fieldCategory": "kyc",
                    "format": "none",
                    "fieldType": "default",
                    "$ref": "#/definitions/simpleType"
                },
                "proofOfException-1": {
                    "bioAttributes": [],
                    "fieldCategory": "evidence",
                    "format": "none",
                    "fieldType": "default",
                    "$ref": "#/definitions/documentType"
                },
                "individualBiometrics": {
                    "bioAttributes": [
                        "leftEye",
                        "rightEye",
                        "rightIndex",
                        "rightLittle",
                        "rightRing",
                        "rightMiddle",
                        "leftIndex",
                        "leftLittle",
                        "leftRing",
                        "leftMiddle",
                        "leftThumb",
                        "rightThumb",
                        "face"
                    ],
                    "fieldCategory": "pvt",
                    "format": "none",
                    "fieldType": "default",
                    "$ref": "#/definitions/biometricsType"
                },
                "addressLine1": {
                    "bioAttributes": [],
                    "validators": [
                        {
                            "validator": "^(?=.{3,50}$).*",
                            "arguments": [],
                            "type": "regex"
                        }
                    ],
                    "fieldCategory": "pvt",
                    "format": "none",
                    "fieldType": "default",
                    "$ref": "#/definitions/simpleType"
                },
                "addressLine2": {
                    "bioAttributes": [],
                    "validators": [
                        {
                            "validator": "^(?=.{3,50}$).*",
                            "arguments": [],
                            "type": "regex"
                        }
                    ],
                    "fieldCategory": "pvt",
                    "format": "none",
                    "fieldType": "default",
                    "$ref": "#/definitions/simpleType"
                },
                "addressLine3": {
                    "bioAttributes": [],
                    "validators": [
                        {
                            "validator": "^(?=.{3,50}$).*",
                            "arguments": [],
                            "type": "regex"
                        }
                    ],
                    "fieldCategory": "pvt",
                    "format": "none",
                    "fieldType": "default",
                    "$ref": "#/definitions/simpleType"
                },
                "email": {
                    "bioAttributes": [],
                    "validators": [
                        {
                            "validator": "^[A-Za-z0-9_\\-]+(\\.[A-Za-z0-9_]+)*@[A-Za-z0-9_-]+(\\.[A-Za-z0-9_]+)*(\\.[a-zA-Z]{2,})$",
                            "arguments": [],
                            "type": "regex"
                        }
                    ],
                    "fieldCategory": "pvt",
                    "format": "none",
                    "type": "string",
                    "fieldType": "default"
                },
                "introducerRID": {
                    "bioAttributes": [],
                    "fieldCategory": "evidence",
                    "format": "none",
                    "type": "string",
                    "fieldType": "default"
                },
                "introducerBiometrics": {
                    "bioAttributes": [
                        "leftEye",
                        "rightEye",
                        "rightIndex",
                        "rightLittle",
                        "rightRing",
                        "rightMiddle",
                        "leftIndex",
                        "leftLittle",
                        "leftRing",
                        "leftMiddle",
                        "leftThumb",
                        "rightThumb",
                        "face"
                    ],
                    "fieldCategory": "pvt",
                    "format": "none",
                    "fieldType": "default",
                    "$ref": "#/definitions/biometricsType"
                },
                "firstName": {
                    "bioAttributes": [],
                    "validators": [
                        {
                            "validator": "^(?=.{2,50}$).*",
                            "arguments": [],
                            "type": "regex"
                        }
                    ],
                    "fieldCategory": "pvt",
                    "format": "none",
                    "fieldType": "default",
                    "$ref": "#/definitions/simpleType"
                },
                "lastName": {
                    "bioAttributes": [],
                    "validators": [
                        {
                            "validator": "^(?=.{2,50}$).*",
                            "arguments": [],
                            "type": "regex"
                        }
                    ],
                    "fieldCategory": "pvt",
                    "format": "none",
                    "fieldType": "default",
                    "$ref": "#/definitions/simpleType"
                },
                "dateOfBirth": {
                    "bioAttributes": [],
                    "validators": [
                        {
                            "validator": "^(1869|18[7-9][0-9]|19[0-9][0-9]|20[0-9][0-9])/([0][1-9]|1[0-2])/([0][1-9]|[1-2][0-9]|3[01])$",
                            "arguments": [],
                            "type": "regex"
                        }
                    ],
                    "fieldCategory": "pvt",
                    "format": "none",
                    "type": "string",
                    "fieldType": "default"
                },
                "individualAuthBiometrics": {
                    "bioAttributes": [
                        "leftEye",
                        "rightEye",
                        "rightIndex",
                        "rightLittle",
                        "rightRing",
                        "rightMiddle",
                        "leftIndex",
                        "leftLittle",
                        "leftRing",
                        "leftMiddle",
                        "leftThumb",
                        "rightThumb",
                        "face"
                    ],
                    "fieldCategory": "pvt",
                    "format": "none",
                    "fieldType": "default",
                    "$ref": "#/definitions/biometricsType"
                },
                "introducerUIN": {
                    "bioAttributes": [],
                    "fieldCategory": "evidence",
                    "format": "none",
                    "type": "string",
                    "fieldType": "default"
                },
                "proofOfIdentity": {
                    "bioAttributes": [],
                    "fieldCategory": "pvt",
                    "format": "none",
                    "fieldType": "default",
                    "$ref": "#/definitions/documentType"
                },
                "IDSchemaVersion": {
                    "bioAttributes": [],
                    "fieldCategory": "none",
                    "format": "none",
                    "type": "number",
                    "fieldType": "default",
                    "minimum": 0
                },
                "proofOfException": {
                    "bioAttributes": [],
                    "fieldCategory": "evidence",
                    "format": "none",
                    "fieldType": "default",
                    "$ref": "#/definitions/documentType"
                },
                "phone": {
                    "bioAttributes": [],
                    "validators": [
                        {
                            "validator": "^[+]*([0-9][0-9]{8,9})$",
                            "arguments": [],
                            "type": "regex"
                        }
                    ],
                    "fieldCategory": "pvt",
                    "format": "none",
                    "type": "string",
                    "fieldType": "default"
                },
                "introducerName": {
                    "bioAttributes": [],
                    "fieldCategory": "evidence",
                    "format": "none",
                    "fieldType": "default",
                    "$ref": "#/definitions/simpleType"
                },
                "proofOfRelationship": {
                    "bioAttributes": [],
                    "fieldCategory": "pvt",
                    "format": "none",
                    "fieldType": "default",
                    "$ref": "#/definitions/documentType"
                },
                "UIN": {
                    "bioAttributes": [],
                    "fieldCategory": "none",
                    "format": "none",
                    "type": "string",
                    "fieldType": "default"
                },
                "preferredLang": {
                    "bioAttributes": [],
                    "fieldCategory": "pvt",
                    "format": "none",
                    "type": "string",
                    "fieldType": "dynamic"
                }
            }
        }
    }
}
'@

$ws.Range("F2").Value = $schemaJson
$ws.Rows.Item(2).RowHeight = 409.6

$ws.Range("H2").Select()
